$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status changes from "Not yet handed off" to "Handed back" for the
#    baf06706... file, on the summary sheet and both language sheets.
# ---------------------------------------------------------------------------
$wsOverview.Range("B2").Value = "Handed back"
$wsOverview.Range("C2").Value = "Handed back"
$wsZh.Range("B2").Value = "Handed back"
$wsDe.Range("B2").Value = "Handed back"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: add "Latest Target File" (E2) / "Latest Handback File" (F2)
#    hyperlinks (re-using the existing Source/Handoff targets) and stamp the
#    "Latest Handback DateTime" (G2).
# ---------------------------------------------------------------------------
$zhLinks = @($wsZh.Hyperlinks)
$zhSourceLink = $zhLinks[0]
$zhHandoffLink = $zhLinks[1]
$zhConfigLink = $zhLinks[2]

$zhSourceAddress = $zhSourceLink.Address
$zhSourceDisplay = $zhSourceLink.TextToDisplay
$zhHandoffAddress = $zhHandoffLink.Address
$zhHandoffDisplay = $zhHandoffLink.TextToDisplay
$zhConfigAddress = $zhConfigLink.Address
$zhConfigDisplay = $zhConfigLink.TextToDisplay

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhSourceAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $zhSourceDisplay) | Out-Null
$wsZh.Range("E2").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhHandoffAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $zhHandoffDisplay) | Out-Null
$wsZh.Range("F2").Style = "HyperLink"

$wsZh.Range("G2").Value = "2016-01-07 11:04:06"

# Re-create the .localization-config hyperlink (row 3) last, so it keeps
# sorting after the two new row-2 hyperlinks (E2, F2).
$zhConfigLink.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhConfigAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $zhConfigDisplay) | Out-Null
$wsZh.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment as zh-cn.
# ---------------------------------------------------------------------------
$deLinks = @($wsDe.Hyperlinks)
$deSourceLink = $deLinks[0]
$deHandoffLink = $deLinks[1]
$deConfigLink = $deLinks[2]

$deSourceAddress = $deSourceLink.Address
$deSourceDisplay = $deSourceLink.TextToDisplay
$deHandoffAddress = $deHandoffLink.Address
$deHandoffDisplay = $deHandoffLink.TextToDisplay
$deConfigAddress = $deConfigLink.Address
$deConfigDisplay = $deConfigLink.TextToDisplay

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deSourceAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $deSourceDisplay) | Out-Null
$wsDe.Range("E2").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deHandoffAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $deHandoffDisplay) | Out-Null
$wsDe.Range("F2").Style = "HyperLink"

$wsDe.Range("G2").Value = "2016-01-07 11:04:20"

$deConfigLink.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deConfigAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $deConfigDisplay) | Out-Null
$wsDe.Range("A3").Style = "HyperLink"
